$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.020.47"
$ws.Range("E2").Value = "  +6.14%  "
$ws.Range("D3").Value = "3.649.66"
$ws.Range("E3").Value = "  +5.82%  "
$ws.Range("D5").Value = "'593.88"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").Value = "'196.00"
$ws.Range("E6").Value = "  +3.53%  "
$ws.Range("D7").Value = "'0.649"
$ws.Range("E7").Value = "  +2.75%  "
$ws.Range("D8").Value = "3.643.25"
$ws.Range("E8").Value = "  +5.86%  "
$ws.Range("D10").Value = "'0.185"
$ws.Range("E10").Value = "  +7.78%  "
$ws.Range("D11").Value = "'0.683"
$ws.Range("E11").Value = "  +5.68%  "
$ws.Range("D12").Value = "'58.37"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("E13").Value = "  +7.11%  "
$ws.Range("D14").Value = "'10.03"
$ws.Range("E14").Value = "  +5.72%  "
$ws.Range("D15").Value = "4.231.35"
$ws.Range("E15").Value = "  +5.60%  "
$ws.Range("D16").Value = "'20.38"
$ws.Range("E16").Value = "  +7.72%  "
$ws.Range("D17").Value = "3.647.26"
$ws.Range("E17").Value = "  +5.63%  "
$ws.Range("D18").Value = "70.952.12"
$ws.Range("E18").Value = "  +5.82%  "
$ws.Range("E19").Value = "  +5.50%  "
$ws.Range("E20").Value = "  +2.48%  "
$ws.Range("E21").Value = "  +4.08%  "
$ws.Range("D22").Value = "'490.29"
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("D23").Value = "'19.42"
$ws.Range("E23").Value = "  +14.45%  "
$ws.Range("D24").Value = "'5.25"
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("D25").Value = "'4.50"
$ws.Range("E25").Value = "  +3.01%  "
$ws.Range("D26").Value = "'91.62"
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("D27").Value = "'3.18"
$ws.Range("E27").Value = "  +6.30%  "
$ws.Range("D28").Value = "'11.49"
$ws.Range("E28").Value = "  +4.91%  "
$ws.Range("D29").Value = "'9.65"
$ws.Range("E29").Value = "  +6.59%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'7.97"
$ws.Range("E30").Value = "  +6.43%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'32.88"
$ws.Range("E31").Value = "  +5.09%  "
$ws.Range("E32").Value = "  +9.99%  "
$ws.Range("D33").Value = "'12.32"
$ws.Range("E33").Value = "  +4.49%  "
$ws.Range("D34").Value = "'617.77"
$ws.Range("E34").Value = "  +2.87%  "
$ws.Range("D35").Value = "'66.59"
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("D36").Value = "'40.41"
$ws.Range("E36").Value = "  +7.86%  "
$ws.Range("D37").Value = "0.0₃0836"
$ws.Range("E37").Value = "  +11.04%  "
$ws.Range("D38").Value = "'0.414"
$ws.Range("E38").Value = "  +5.84%  "
$ws.Range("E39").Value = "  +1.21%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  +2.52%  "
$ws.Range("D42").Value = "3.337.19"
$ws.Range("E42").Value = "  +3.66%  "
$ws.Range("D43").Value = "'3.24"
$ws.Range("E43").Value = "  +16.16%  "
$ws.Range("D44").Value = "'3.18"
$ws.Range("E44").Value = "  +8.42%  "
$ws.Range("D45").Value = "'2.86"
$ws.Range("E45").Value = "  +10.25%  "
$ws.Range("D46").Value = "'0.0460"
$ws.Range("E46").Value = "  +6.60%  "
$ws.Range("D47").Value = "'9.67"
$ws.Range("E47").Value = "  +11.05%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.140"
$ws.Range("E48").Value = "  +3.45%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "'3.32"
$ws.Range("E49").Value = "  +2.48%  "
$ws.Range("D50").Value = "'3.25"
$ws.Range("E50").Value = "  +1.48%  "
$ws.Range("E51").Value = "  -0.15%  "
